# Fruta / hortaliza, semanal
# Insert a new week's worth of data (3 rows) at the top of the
# "Terminal La Palmera de La Serena - Plátano" price block (rows 349-400),
# pushing the existing rows down by one group of 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 349 - this shifts rows 349:400 down to 352:403
# and Excel copies the row-349 formatting (incl. the date NumberFormat on
# column D) into the freshly inserted rows, same as native Excel behaviour.
$ws.Range("A349:A351").EntireRow.Insert()

# Row 349 - "Pintón"
$ws.Cells.Item(349, 1).Value = 8
$ws.Cells.Item(349, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(349, 3).Value = "Coquimbo"
$ws.Cells.Item(349, 4).Value = 44522
$ws.Cells.Item(349, 5).Value = 4
$ws.Cells.Item(349, 6).Value = "Fruta"
$ws.Cells.Item(349, 7).Value = 100108
$ws.Cells.Item(349, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(349, 9).Value = 100108006
$ws.Cells.Item(349, 10).Value = "Plátano"
$ws.Cells.Item(349, 11).Value = "Sin especificar"
$ws.Cells.Item(349, 12).Value = "Pintón"
$ws.Cells.Item(349, 13).Value = 80
$ws.Cells.Item(349, 14).Value = 22000
$ws.Cells.Item(349, 15).Value = 22000
$ws.Cells.Item(349, 16).Value = 22000
$ws.Cells.Item(349, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(349, 18).Value = "Ecuador"
$ws.Cells.Item(349, 19).Value = 1100
$ws.Cells.Item(349, 20).Value = 20

# Row 350 - "Primera Maduro"
$ws.Cells.Item(350, 1).Value = 8
$ws.Cells.Item(350, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(350, 3).Value = "Coquimbo"
$ws.Cells.Item(350, 4).Value = 44522
$ws.Cells.Item(350, 5).Value = 4
$ws.Cells.Item(350, 6).Value = "Fruta"
$ws.Cells.Item(350, 7).Value = 100108
$ws.Cells.Item(350, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(350, 9).Value = 100108006
$ws.Cells.Item(350, 10).Value = "Plátano"
$ws.Cells.Item(350, 11).Value = "Sin especificar"
$ws.Cells.Item(350, 12).Value = "Primera Maduro"
$ws.Cells.Item(350, 13).Value = 120
$ws.Cells.Item(350, 14).Value = 23000
$ws.Cells.Item(350, 15).Value = 23000
$ws.Cells.Item(350, 16).Value = 23000
$ws.Cells.Item(350, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(350, 18).Value = "Ecuador"
$ws.Cells.Item(350, 19).Value = 1150
$ws.Cells.Item(350, 20).Value = 20

# Row 351 - "Primera Pintón"
$ws.Cells.Item(351, 1).Value = 8
$ws.Cells.Item(351, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(351, 3).Value = "Coquimbo"
$ws.Cells.Item(351, 4).Value = 44522
$ws.Cells.Item(351, 5).Value = 4
$ws.Cells.Item(351, 6).Value = "Fruta"
$ws.Cells.Item(351, 7).Value = 100108
$ws.Cells.Item(351, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(351, 9).Value = 100108006
$ws.Cells.Item(351, 10).Value = "Plátano"
$ws.Cells.Item(351, 11).Value = "Sin especificar"
$ws.Cells.Item(351, 12).Value = "Primera Pintón"
$ws.Cells.Item(351, 13).Value = 120
$ws.Cells.Item(351, 14).Value = 24000
$ws.Cells.Item(351, 15).Value = 24000
$ws.Cells.Item(351, 16).Value = 24000
$ws.Cells.Item(351, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(351, 18).Value = "Ecuador"
$ws.Cells.Item(351, 19).Value = 1200
$ws.Cells.Item(351, 20).Value = 20
